# calc correct avg exp times
# Update Avg_Agent_Step_Time (G), Avg_Experiment_Time (H),
# Std_Agent_Step_Time (M) and Std_Experiment_Time (N) columns
# for rows 2-13 with recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.470898529999998
$ws.Range("H2").Value = 481.82970353
$ws.Range("M2").Value = 1.038744466823721
$ws.Range("N2").Value = 87.85142470676432

$ws.Range("G3").Value = 7.929341020000001
$ws.Range("H3").Value = 687.9941678200001
$ws.Range("M3").Value = 0.7531221812524539
$ws.Range("N3").Value = 116.4694584567193

$ws.Range("G4").Value = 3.31464951
$ws.Range("H4").Value = 101.15000589
$ws.Range("M4").Value = 0.4536641048678035
$ws.Range("N4").Value = 25.54897390697773

$ws.Range("G5").Value = 3.24284849
$ws.Range("H5").Value = 151.69987375
$ws.Range("M5").Value = 0.3429550758920183
$ws.Range("N5").Value = 35.14484673333508

$ws.Range("G6").Value = 1.10917767
$ws.Range("H6").Value = 17.17963813
$ws.Range("M6").Value = 0.2534476437025396
$ws.Range("N6").Value = 6.769580720283482

$ws.Range("G7").Value = 1.20927583
$ws.Range("H7").Value = 29.6374162
$ws.Range("M7").Value = 0.1999060830073323
$ws.Range("N7").Value = 9.679015808747899

$ws.Range("G8").Value = 0.54786952
$ws.Range("H8").Value = 5.67933464
$ws.Range("M8").Value = 0.1726342034945155
$ws.Range("N8").Value = 2.700303271253178

$ws.Range("G9").Value = 0.6083432400000001
$ws.Range("H9").Value = 10.57761577
$ws.Range("M9").Value = 0.1404447433399849
$ws.Range("N9").Value = 4.748946264075736

$ws.Range("G10").Value = 0.29505848
$ws.Range("H10").Value = 2.354936260000001
$ws.Range("M10").Value = 0.1194965853319483
$ws.Range("N10").Value = 1.390704615183722

$ws.Range("G11").Value = 0.33874006
$ws.Range("H11").Value = 4.62158223
$ws.Range("M11").Value = 0.09368018421475756
$ws.Range("N11").Value = 2.390849924646737

$ws.Range("G12").Value = 0.18407432
$ws.Range("H12").Value = 1.22188418
$ws.Range("M12").Value = 0.08572415951680402
$ws.Range("N12").Value = 0.8353332643361175

$ws.Range("G13").Value = 0.2046409
$ws.Range("H13").Value = 2.34155125
$ws.Range("M13").Value = 0.07017460458968985
$ws.Range("N13").Value = 1.477540619396944
